$wb = $excel.ActiveWorkbook

$prices1 = @{
    2 = 40.3
    3 = 33.65
    4 = 49.32
    5 = 49.15
    6 = 28.6
    7 = 41.46
    8 = 25.6
    9 = 30.4
    10 = 30.94
    11 = 109.9
    12 = 140.31
    13 = 179.42
    14 = 0.6919999999999999
    15 = 9.210000000000001
    16 = 26.75
    17 = 26.71
    18 = 14.29
    19 = 35.18
    20 = 32.17
    21 = 21.78
    22 = 83.69
    23 = 3.499
    24 = 4.648
}

$prices2 = @{
    2 = 40.3
    3 = 33.65
    4 = 49.32
    5 = 49.15
    6 = 28.6
    7 = 41.46
    8 = 25.6
    9 = 30.4
    10 = 30.94
    11 = 109.9
    12 = 179.42
    13 = 140.31
    14 = 230.73
    15 = 0.6919999999999999
    16 = 37.95
    17 = 9.210000000000001
    18 = 26.75
    19 = 26.71
    20 = 14.29
    21 = 35.18
    22 = 32.17
    23 = 21.78
    24 = 83.69
    25 = 3.499
    26 = 0.865
    27 = 4.648
}

$ws1 = $wb.Worksheets.Item("个人持仓")
foreach ($row in $prices1.Keys) {
    $ws1.Range("C$row").Value = $prices1[$row]
}

$ws2 = $wb.Worksheets.Item("家庭持仓")
foreach ($row in $prices2.Keys) {
    $ws2.Range("C$row").Value = $prices2[$row]
}
